$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-DTCP v1.2")

# Row 13 (VJ0603D471KXXAJ): designator list trimmed down and quantity reduced to match.
$ws.Range("C13").Value = "C18, C43, C45, C59"
$ws.Range("F13").Value = 4

# Row 30 previously held the "1812PS-333_R_" group (L5, L6, L10). Replace it in place
# with the new LQW2BHN68NJ03L part covering L5, L6, L7.
$ws.Range("A30").Value = "LQW2BHN68NJ03L"
$ws.Range("B30").Value = "Wire Wound RF Inductor 68nH ±5% 460mA 0.23" + [char]0x03A9 + " 0805 (2015)"
$ws.Range("C30").Value = "L5, L6, L7"
$ws.Range("D30").Value = "FP-LQW2BHN_03-MFG"
$ws.Range("E30").Value = "CMP-06042-008671-1"
$ws.Range("F30").Value = 3

# Row 31 previously held "LQG18HHR10J00D" (L7, L8, L9, L11, L12). Replace it in place
# with the remaining "1812PS-333_R_" part, now only covering L10.
$ws.Range("A31").Value = "1812PS-333_R_"
$ws.Range("B31").Value = "1812PS-333_R_ Coilcraft"
$ws.Range("C31").Value = "L10"
$ws.Range("D31").Value = "FP-1812PS-MFG"
$ws.Range("E31").Value = "CMP-11268-000188-1"
$ws.Range("F31").Value = 1

# Row 32 ("LQG18HH82NJ00D", L13) is no longer needed now that its inductors were folded
# into the rows above; delete the whole row so everything below shifts up by one.
$ws.Rows("32").Delete()
